$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 76, shifting existing rows 76..122 down to 77..123
$ws.Rows.Item(76).Insert()

# Populate the new row 76 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,N,Q,R are constant across every data row in this sheet.
$ws.Cells.Item(76, 1).Value2 = 7
$ws.Cells.Item(76, 2).Value2 = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(76, 3).Value2 = "Ñuble"
$ws.Cells.Item(76, 4).Value2 = 45001
$ws.Cells.Item(76, 5).Value2 = 16
$ws.Cells.Item(76, 6).Value2 = 100112030
$ws.Cells.Item(76, 7).Value2 = "Poroto granado"
$ws.Cells.Item(76, 8).Value2 = "Sin especificar"
$ws.Cells.Item(76, 9).Value2 = "Primera"
$ws.Cells.Item(76, 10).Value2 = 30
$ws.Cells.Item(76, 11).Value2 = 30000
$ws.Cells.Item(76, 12).Value2 = 30000
$ws.Cells.Item(76, 13).Value2 = 30000
$ws.Cells.Item(76, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(76, 15).Value2 = "Provincia de Diguillín"
$ws.Cells.Item(76, 16).Value2 = 1200
$ws.Cells.Item(76, 17).Value2 = 25
$ws.Cells.Item(76, 18).Value2 = "Hortaliza"
